# Update B15 from "Job title, Position" to "Job title"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B15").Value = "Job title"

# Add 3 new rows (67-69) at the bottom, copying the format of row 66
# (A column style = s4, B column style = s6) then set their text values.

$ws.Range("A66:B66").Copy($ws.Range("A67:B67"))
$ws.Range("A67").Value = "Quá trình công tác"
$ws.Range("B67").Value = "Working Progress"
$ws.Rows.Item(67).RowHeight = 15.75

$ws.Range("A66:B66").Copy($ws.Range("A68:B68"))
$ws.Range("A68").Value = "Quá trình tập sự"
$ws.Range("B68").Value = "Probation Progress"
$ws.Rows.Item(68).RowHeight = 15.75

$ws.Range("A66:B66").Copy($ws.Range("A69:B69"))
$ws.Range("A69").Value = "Diễn biến lương"
$ws.Range("B69").Value = "Wage Changes"
$ws.Rows.Item(69).RowHeight = 15.75

# Update selection / view to focus near the newly added rows
$ws.Activate()
$ws.Range("B70").Select()
